# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51) to
# reflect freshly scraped figures. Column D ("Price") and column E
# ("Volume(1h)") values are refreshed for every row except row 4
# (TetherUSD, unchanged), and rows 30/31 plus 42/43 additionally swap
# which coin occupies which row (their Coin name + Link also change).
#
# NOTE: every Price/Volume cell in this sheet is stored as literal TEXT
# (e.g. "45.804.98" uses dots as thousands separators, so it is not a
# valid number anyway; percentages keep their padding spaces). Excel's
# Range.Value setter auto-infers the literal's type, so a plain-looking
# numeric string such as "312.28" or "1.00" would silently become the
# *number* 312.28 / 1 (losing the trailing zero) unless the cell is
# pre-formatted as Text ("@"). We therefore flip NumberFormat to "@" on
# just the cells whose new value would otherwise be misread as a number
# before writing them; cells whose new text is not a valid numeric
# literal (e.g. "45.804.98", with two dots) are left alone since Excel
# already keeps those as text natively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: force text format on D-column cells whose new values would
# otherwise be auto-coerced to numbers by Excel's type inference ---
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16", "D17", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D35", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Step 2: write the updated values ---
$ws.Range("D2").Value = '45.804.98'
$ws.Range("E2").Value = '  -3.11%  '
$ws.Range("D3").Value = '2.680.67'
$ws.Range("E3").Value = '  +1.50%  '
$ws.Range("D5").Value = '312.28'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '98.65'
$ws.Range("E6").Value = '  -6.06%  '
$ws.Range("D7").Value = '0.599'
$ws.Range("E7").Value = '  -2.58%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '0.584'
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").Value = '38.39'
$ws.Range("E10").Value = '  -3.73%  '
$ws.Range("D11").Value = '0.0849'
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = '8.17'
$ws.Range("E12").Value = '  -3.47%  '
$ws.Range("D13").Value = '3.088.15'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("E14").Value = '  +0.81%  '
$ws.Range("D15").Value = '2.674.68'
$ws.Range("E15").Value = '  +1.03%  '
$ws.Range("D16").Value = '0.934'
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").Value = '15.21'
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").Value = '45.805.23'
$ws.Range("E18").Value = '  -4.01%  '
$ws.Range("E19").Value = '  -1.47%  '
$ws.Range("D20").Value = '6.87'
$ws.Range("E20").Value = '  +0.77%  '
$ws.Range("D21").Value = '12.90'
$ws.Range("E21").Value = '  -3.29%  '
$ws.Range("D22").Value = '75.38'
$ws.Range("E22").Value = '  +3.35%  '
$ws.Range("D23").Value = '283.67'
$ws.Range("E23").Value = '  +4.10%  '
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '31.45'
$ws.Range("E25").Value = '  +2.85%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '4.06'
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '10.59'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '38.36'
$ws.Range("E30").Value = '  -4.80%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '2.17'
$ws.Range("E31").Value = '  -6.36%  '
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("D35").Value = '155.30'
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").Value = '0.0843'
$ws.Range("E36").Value = '  -1.74%  '
$ws.Range("E37").Value = '  -1.74%  '
$ws.Range("E38").Value = '  -2.23%  '
$ws.Range("D39").Value = '26.06'
$ws.Range("E39").Value = '  +10.20%  '
$ws.Range("D40").Value = '0.125'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '16.27'
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("B42").Value = 'NEARProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D42").Value = '3.62'
$ws.Range("E42").Value = '  -3.94%  '
$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").Value = '0.0329'
$ws.Range("E43").Value = '  -1.89%  '
$ws.Range("D44").Value = '3.98'
$ws.Range("E44").Value = '  -7.38%  '
$ws.Range("D45").Value = '2.152.06'
$ws.Range("E45").Value = '  -1.47%  '
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("D47").Value = '94.51'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '9.36'
$ws.Range("E48").Value = '  -6.66%  '
$ws.Range("D49").Value = '112.20'
$ws.Range("E49").Value = '  -2.13%  '
$ws.Range("D50").Value = '2.939.15'
$ws.Range("E50").Value = '  +1.46%  '
$ws.Range("D51").Value = '0.201'
$ws.Range("E51").Value = '  -2.12%  '
